$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '37.280.71'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  +1.98%  '
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '2.024.96'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  +2.65%  '
$ws.Range("E4").Value = '  -0.10%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '248.20'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +1.71%  '
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '0.621'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  -0.93%  '
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '58.03'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  -3.52%  '
$ws.Range("E8").Value = '  -0.05%  '
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.389'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  +2.88%  '
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '0.0805'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +2.05%  '
$ws.Range("E11").Value = '  -0.39%  '
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '14.94'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  +4.93%  '
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '2.321.40'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  +2.62%  '
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '0.833'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  -1.28%  '
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '21.42'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  -1.28%  '
$ws.Range("E16").Value = '  +0.40%  '
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '2.024.71'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  +2.90%  '
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '37.219.44'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  +1.83%  '
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '69.95'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  +0.19%  '
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '0.0₃0857'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  +0.14%  '
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '5.23'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  +2.87%  '
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '228.58'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  -0.43%  '
$ws.Range("E23").Value = '  +0.07%  '
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '2.54'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  +4.60%  '
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '2.35'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  -0.62%  '
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '9.19'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  +0.60%  '
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '163.77'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  +0.75%  '
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '0.138'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  -5.30%  '
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '19.90'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  +2.69%  '
$ws.Range("E30").Value = '  +2.43%  '
$ws.Range("E31").Value = '  +0.03%  '
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '4.77'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  -0.48%  '
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '0.0667'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  +8.44%  '
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '4.56'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  +0.98%  '
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '2.49'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  +8.82%  '
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '3.53'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  +5.77%  '
$ws.Range("E37").Value = '  -0.09%  '
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '1.82'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  +2.47%  '
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '5.36'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  -0.61%  '
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '2.99'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  +2.82%  '
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '0.0970'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  -0.16%  '
$ws.Range("E42").Value = '  +3.49%  '
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '1.17'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  -0.07%  '
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '16.44'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  +3.13%  '
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '1.398.21'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  +2.39%  '
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '91.00'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  +2.06%  '
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '7.48'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  +4.24%  '
$ws.Range("E48").Value = '  +1.93%  '
$ws.Range("E49").Value = '  +11.89%  '
$ws.Range("E50").Value = '  +1.47%  '
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '2.213.28'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  +2.65%  '
